$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Status table (Table 1): two status cells change text.
#    Row 3 / Col 2 (Brukerregistrering row): "Plan" -> longer status text.
#    Row 8 / Col 2 (Backend/roller row): "Plan" -> "Påbegynt".
#    Scope the Find to each individual cell's Range and use wdReplaceOne (1)
#    so only the text inside that specific cell is touched (there is also an
#    unrelated "Plan" substring inside "Menu Planner" on page 1).
# ---------------------------------------------------------------------------
$table1 = $d.Tables.Item(1)

$statusCellBrukerreg = $table1.Cell(3, 2)
$statusCellBrukerreg.Range.Find.Execute("Plan", $true, $true, $false, $false, $false, $true, 1, $false, `
    "Delvis gjennomført. Mangler flere innloggingstjenester, og trenger kvalitetstesting.", 1)

$statusCellBackend = $table1.Cell(8, 2)
$statusCellBackend.Range.Find.Execute("Plan", $true, $true, $false, $false, $false, $true, 1, $false, `
    "Påbegynt", 1)

# ---------------------------------------------------------------------------
# Helper pattern used below: Paragraph.Range.Text always carries a trailing
# paragraph-mark (chr 13) so exact "-eq" comparisons never hit; strip it
# before comparing / before feeding the text back into Find as the
# "replace with" value (re-using the exact same text merges the split runs
# into a single run without altering the visible content).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2) "Publisering av ingredienser" bullet: merge the two runs
#    "Ha et godt bilde:" + " Minst 1000x350, ikke oppskalert fra et lite
#    bilde, «profesjonell» fotokvalitet." into a single run with identical
#    combined text. There is an earlier, unrelated paragraph (under
#    "Publisering av oppskrifter") with the very same final text that must
#    stay untouched, so identify the right one via its paragraph index.
# ---------------------------------------------------------------------------
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text.Contains("fotokvalitet") -and $i -gt 290) {
        $clean = $text.TrimEnd([char]13)
        $para.Range.Find.Execute($clean, $true, $true, $false, $false, $false, $true, 1, $false, $clean, 1)
        break
    }
}

# ---------------------------------------------------------------------------
# 3) "Publisering av næringsstoffer" intro: merge the three runs
#    "E" + "t næringsstoff" + " må oppfylle følgende krav for å kunne
#    publiseres:" into a single run.
# ---------------------------------------------------------------------------
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text.Contains("næringsstoff må oppfylle")) {
        $clean = $text.TrimEnd([char]13)
        $para.Range.Find.Execute($clean, $true, $true, $false, $false, $false, $true, 1, $false, $clean, 1)
        break
    }
}

# ---------------------------------------------------------------------------
# 4) The bullet right after it: merge the three runs
#    "Ha et godt bilde" + " eller illustrasjon" + ": Minst 1000x350, ikke
#    oppskalert fra et lite bilde." into a single run.
# ---------------------------------------------------------------------------
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text.Contains("eller illustrasjon")) {
        $clean = $text.TrimEnd([char]13)
        $para.Range.Find.Execute($clean, $true, $true, $false, $false, $false, $true, 1, $false, $clean, 1)
        break
    }
}
